$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 249; this pushes the existing
# rows 249..344 down to 250..345 (and bumps the sheet dimension to
# A1:R345), exactly like Excel's own "Insert Sheet Rows" command.
$ws.Rows("249:249").Insert()

# Populate the newly inserted row 249 with its data (same categorical
# fields as the row that used to occupy 249, new measured values).
$ws.Range("A249").Value = 6
$ws.Range("B249").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C249").Value = "Metropolitana"
$ws.Range("D249").Value = 44460
$ws.Range("E249").Value = 13
$ws.Range("F249").Value = 100112003
$ws.Range("G249").Value = "Ajo"
$ws.Range("H249").Value = "Chino"
$ws.Range("I249").Value = "Primera"
$ws.Range("J249").Value = 1600
$ws.Range("K249").Value = 14000
$ws.Range("L249").Value = 14500
$ws.Range("M249").Value = 14312
$ws.Range("N249").Value = "$/caja 10 kilos"
$ws.Range("O249").Value = "China"
$ws.Range("P249").Value = 1431
$ws.Range("Q249").Value = 10
$ws.Range("R249").Value = "Hortaliza"
